# Trade #10 closed at 2026-02-17 20:03:00 - unknown UNKNOWN +0.000%

$wb = $excel.ActiveWorkbook

# --- Sheet: Summary ---
$ws = $wb.Worksheets.Item("Summary")
$ws.Range("B3").Value = 1399.67
$ws.Range("B4").Value = -0.33
$ws.Range("B5").Value = -0.66
$ws.Range("B6").Value = 10
$ws.Range("B8").Value = 7
$ws.Range("B9").Value = 30

# --- Sheet: Strategy Status ---
$ws = $wb.Worksheets.Item("Strategy Status")
$ws.Range("C5").Value = 99.67
$ws.Range("D5").Value = 10
$ws.Range("E5").Value = -0.33
$ws.Range("F5").Value = -0.33
$ws.Range("G5").Value = 30

# --- Sheet: All Trades ---
$ws = $wb.Worksheets.Item("All Trades")
$ws.Range("A11").Value = 10
$ws.Range("B11").NumberFormat = "@"
$ws.Range("B11").Value = "2026-02-17"
$ws.Range("B11").ClearFormats()
$ws.Range("C11").NumberFormat = "@"
$ws.Range("C11").Value = "20:02:54"
$ws.Range("C11").ClearFormats()
$ws.Range("D11").Value = "MarketMaking"
$ws.Range("E11").Value = "DOWN"
$ws.Range("F11").Value = 0.309278
$ws.Range("G11").Value = 0.27
$ws.Range("H11").Value = "CLOSED"
$ws.Range("I11").Value = -12.7
$ws.Range("J11").Value = -0.04
$ws.Range("K11").Value = 99.67
$ws.Range("L11").Value = 0
$ws.Range("M11").Value = 0
$ws.Range("N11").Value = 0.6
$ws.Range("O11").Value = "Normal spread capture: 19600 bps"
$ws.Range("P11").Value = "early_exit"
$ws.Range("Q11").Value = 0.13

# --- Sheet: MarketMaking ---
$ws = $wb.Worksheets.Item("MarketMaking")
$ws.Range("A11").Value = 10
$ws.Range("B11").NumberFormat = "@"
$ws.Range("B11").Value = "2026-02-17"
$ws.Range("B11").ClearFormats()
$ws.Range("C11").NumberFormat = "@"
$ws.Range("C11").Value = "20:02:54"
$ws.Range("C11").ClearFormats()
$ws.Range("D11").Value = "MarketMaking"
$ws.Range("E11").Value = "DOWN"
$ws.Range("F11").Value = 0.309278
$ws.Range("G11").Value = 0.27
$ws.Range("H11").Value = "CLOSED"
$ws.Range("I11").Value = -12.7
$ws.Range("J11").Value = -0.04
$ws.Range("K11").Value = 99.67
$ws.Range("L11").Value = 0
$ws.Range("M11").Value = 0
$ws.Range("N11").Value = 0.6
$ws.Range("O11").Value = "Normal spread capture: 19600 bps"
$ws.Range("P11").Value = "early_exit"
$ws.Range("Q11").Value = 0.13
